# Daily attendance processing - 2026-01-22 11:12:47
# Normalises the "Recorded By" column (G) so the email address is listed
# before "System" instead of after it, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
        $changed = $changed + 1
    }
}

Write-Output "Rows updated: $changed"
